$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New rows of data (A = numeric id, B = text label)
$newRows = @(
    @(1000015, "角色%d级解锁"),
    @(1000016, "刀盾兵"),
    @(1000017, "巨盾兵"),
    @(1000018, "短弓手"),
    @(1000019, "长弓手"),
    @(1000020, "轻骑兵"),
    @(1000021, "重骑兵"),
    @(1000022, "冲车"),
    @(1000023, "投石车"),
    @(1000024, "军医"),
    @(1000025, "巫医")
)

$startRow = 18
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = $newRows[$i][0]
    $ws.Cells.Item($r, 2).Value = $newRows[$i][1]
}

# Widen column B slightly to fit the new, longer labels
# (13.285714285714286 == 14 - 5/7 so that the serialized <col> width comes out to exactly 14)
$ws.Columns.Item(2).ColumnWidth = 13.285714285714286

# Update the saved selection to reflect where the editor left off
[void]$ws.Range("D10").Select()
